# Bot5 GUI --version 1.00 ultimas correcciones
# Sets every "NO" flag in the AGENCIAS sheet (column F, rows 4-38, skipping
# the blank separator row 21 and the row that was already "SI") to "SI",
# which drops the now-unused "NO" shared string and re-scopes the F4:F38
# list validation around the still-blank F21 cell. Also moves the active
# selection to D13, matching the author's last edit location.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AGENCIAS")

# Flip every "NO" -> "SI" in column F (rows 4-38), leaving the blank
# separator row (21) and any cell already marked "SI" untouched.
for ($r = 4; $r -le 38; $r++) {
    if ($r -eq 21) { continue }
    $cell = $ws.Range("F" + $r)
    if ($cell.Value2() -eq "NO") {
        $cell.Value = "SI"
    }
}

# F21 has no flag value (separator row) - dropping its validation rule
# splits the F4:F38 sqref into F4:F20 / F22:F38, same as the source edit.
$ws.Range("F21").Validation.Delete()

# Move / record the active selection like the author's last interaction.
$ws.Activate()
$ws.Range("D13").Select()
